$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.652.50'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.288.02'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'112.47"
$ws.Range("E5").Value = '  +16.77%  '
$ws.Range("D6").Value = "'267.93"
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").Value = "'0.612"
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").Value = "'47.15"
$ws.Range("E10").Value = '  +4.31%  '
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("E12").Value = '  +8.75%  '
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("D14").Value = "'15.51"
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").Value = '2.627.56'
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").Value = "'0.841"
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").Value = '2.286.02'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = '43.496.56'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = "'6.54"
$ws.Range("E20").Value = '  +5.72%  '
$ws.Range("D21").Value = "'72.19"
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").Value = "'2.50"
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").Value = "'232.52"
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = "'9.46"
$ws.Range("E24").Value = '  +3.18%  '
$ws.Range("D25").Value = "'2.81"
$ws.Range("E25").Value = '  +12.91%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = "'11.35"
$ws.Range("E27").Value = '  +1.37%  '
$ws.Range("D28").Value = "'42.83"
$ws.Range("E28").Value = '  +5.79%  '
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D31").Value = "'176.09"
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").Value = "'21.62"
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("E33").Value = '  +4.45%  '
$ws.Range("D34").Value = "'5.47"
$ws.Range("E34").Value = '  +1.86%  '
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("E36").Value = '  +6.84%  '
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = "'3.80"
$ws.Range("E39").Value = '  +12.76%  '
$ws.Range("E40").Value = '  +3.57%  '
$ws.Range("D41").Value = "'73.77"
$ws.Range("E41").Value = '  +15.35%  '
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("D43").Value = "'13.35"
$ws.Range("E43").Value = '  +9.62%  '
$ws.Range("E44").Value = '  +5.87%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = "'5.94"
$ws.Range("E46").Value = '  +13.53%  '
$ws.Range("D47").Value = "'8.73"
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("E48").Value = '  -2.31%  '
$ws.Range("D49").Value = "'102.04"
$ws.Range("E49").Value = '  +3.87%  '
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("D51").Value = "'0.446"
$ws.Range("E51").Value = '  +4.26%  '

# Reset style on cells that needed a quote-prefix to stay text,
# so the style index matches the original (no explicit style).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
